# Daily attendance processing - 2025-11-16 08:50:31
# Reorders the "Recorded By" (column G) values so that any "System"/"system"
# tokens in the comma-separated list are moved to the front, preserving the
# relative order of the System tokens and of the remaining tokens.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Column G is the "Recorded By" column (see header in row 1)
$col = 7

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ","
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        $systemParts = @()
        $otherParts = @()
        foreach ($p in $trimmed) {
            if ($p.ToLower() -eq "system") {
                $systemParts += $p
            } else {
                $otherParts += $p
            }
        }

        $newParts = $systemParts + $otherParts
        $newValue = [string]::Join(", ", $newParts)

        if ($newValue -ne $value) {
            $cell.Value = $newValue
        }
    }
}
